$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Slide 1 - title text tweak: "... planning" -> "... planning/management"
# ---------------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "Time management and project planning/management"

# ---------------------------------------------------------------------------
# 2. Insert a brand-new "Title and Content" slide right before the existing
#    last ("Thankyou") slide. Because it is inserted *before* that slide,
#    the old slide (with its own id) simply shifts down to the last
#    position - exactly mirroring the id re-shuffle seen in the diff
#    (new slide becomes id 285 / position 25, old "Thankyou" slide - still
#    id 278 - becomes position 26).
# ---------------------------------------------------------------------------
$contentLayout = $p.SlideMaster.CustomLayouts.Item(2)   # "Title and Content"
$newSlide = $p.Slides.AddSlide(25, $contentLayout)

# Title placeholder
$titleTr = $newSlide.Shapes.Item(1).TextFrame.TextRange
$titleTr.Text = "Good project manager attributes"
$titleTr.LanguageID = "en-AU"

# Content placeholder - bullet list. Build it paragraph by paragraph with
# InsertAfter so every new paragraph inherits the language/formatting of
# the one before it instead of resetting to the engine defaults.
$body = $newSlide.Shapes.Item(2).TextFrame.TextRange
$body.Text = "Follow intent rather than letter of the plan "
$body.Font.Size = 28
$body.LanguageID = "en-AU"

$body.InsertAfter("`rServant leader") | Out-Null
$body.InsertAfter("`rStrike a good balance between being engaged / supportive and micromanaging") | Out-Null
$body.InsertAfter("`rSet expectations early and often") | Out-Null
$body.InsertAfter("`rMinimise governance overhead") | Out-Null
$body.InsertAfter("`rStart small and inspect and adapt") | Out-Null
$body.InsertAfter("`rCreate an autonomy, mastery, purpose environment") | Out-Null

# 3rd paragraph ("Strike a good balance ...") is a level-2 bullet at a
# smaller font size.
$para3 = $body.Paragraphs(3, 1)
$para3.IndentLevel = 2
$para3.Font.Size = 24

# Slide-number placeholder: copy the one from an existing "Title and
# Content" slide so the field keeps the deck's shared field GUID / styling
# instead of minting a brand-new placeholder from scratch.
$templateSlideNum = $p.Slides.Item(2).Shapes.Item(3)
$templateSlideNum.Copy()
$newSlide.Shapes.Paste() | Out-Null

Write-Output "done"
